$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 221.53847
$ws.Range("I2").Value = 190
$ws.Range("J2").Value = 272
$ws.Range("K2").Value = 190
$ws.Range("L2").Value = 272
$ws.Range("M2").Value = -77
$ws.Range("N2").Value = -498

$ws.Range("H9").Value = 123.333336
$ws.Range("I9").Value = 50
$ws.Range("J9").Value = 138
$ws.Range("K9").Value = 50
$ws.Range("L9").Value = 138
$ws.Range("M9").Value = 119
$ws.Range("N9").Value = -476

$ws.Range("H112").Value = 1435.5416
$ws.Range("I112").Value = 720
$ws.Range("J112").Value = 1518.7441
$ws.Range("K112").Value = 2160
$ws.Range("L112").Value = 4556.2323
$ws.Range("M112").Value = -1052
$ws.Range("N112").Value = -6772.2323

$ws.Range("H129").Value = 704.8200000000001
$ws.Range("I129").Value = 412.5
$ws.Range("J129").Value = 717
$ws.Range("K129").Value = 1237.5
$ws.Range("L129").Value = 2151
$ws.Range("M129").Value = 3762.5
$ws.Range("N129").Value = -12151

$ws.Range("H137").Value = 9256232
$ws.Range("I137").Value = 13758208
$ws.Range("J137").Value = 2169.9443
$ws.Range("K137").Value = 41274624
$ws.Range("L137").Value = 6509.8329
$ws.Range("M137").Value = -41272074
$ws.Range("N137").Value = -11609.8329

$ws.Range("H138").Value = 2499.29
$ws.Range("I138").Value = 1061.5769
$ws.Range("J138").Value = 3004.4324
$ws.Range("K138").Value = 3184.7307
$ws.Range("L138").Value = 9013.297200000001
$ws.Range("M138").Value = 1955.2693
$ws.Range("N138").Value = -19293.2972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2372.6667
$ws.Range("I61").Value = 2042.2727
$ws.Range("J61").Value = 6007
$ws.Range("K61").Value = 2042.2727
$ws.Range("L61").Value = 6007
$ws.Range("M61").Value = -1830.2727
$ws.Range("N61").Value = -6431

$ws.Range("H74").Value = 1149.16
$ws.Range("I74").Value = 944.5238000000001
$ws.Range("J74").Value = 2223.5
$ws.Range("K74").Value = 944.5238000000001
$ws.Range("L74").Value = 2223.5
$ws.Range("M74").Value = -70.52380000000005
$ws.Range("N74").Value = -3971.5

$ws.Range("H77").Value = 1149.16
$ws.Range("I77").Value = 944.5238000000001
$ws.Range("J77").Value = 2223.5
$ws.Range("K77").Value = 4722.619000000001
$ws.Range("L77").Value = 11117.5
$ws.Range("M77").Value = -354.6190000000006
$ws.Range("N77").Value = -19853.5

$ws.Range("H102").Value = 3545.3845
$ws.Range("I102").Value = 2544.2856
$ws.Range("J102").Value = 4713.3335
$ws.Range("K102").Value = 2544.2856
$ws.Range("L102").Value = 4713.3335
$ws.Range("M102").Value = -922.2856000000002
$ws.Range("N102").Value = -7957.3335

$ws.Range("H122").Value = 2705.2856
$ws.Range("I122").Value = 2984.6
$ws.Range("J122").Value = 2007
$ws.Range("K122").Value = 8953.799999999999
$ws.Range("L122").Value = 6021
$ws.Range("M122").Value = -6503.799999999999
$ws.Range("N122").Value = -10921

$ws.Range("H136").Value = 2372.6667
$ws.Range("I136").Value = 2042.2727
$ws.Range("J136").Value = 6007
$ws.Range("K136").Value = 6126.8181
$ws.Range("L136").Value = 18021
$ws.Range("M136").Value = -3576.8181
$ws.Range("N136").Value = -23121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1549.1305
$ws.Range("I105").Value = 1122.8334
$ws.Range("J105").Value = 2014.1818
$ws.Range("K105").Value = 1122.8334
$ws.Range("L105").Value = 2014.1818
$ws.Range("M105").Value = 624.1666
$ws.Range("N105").Value = -5508.1818

$ws.Range("H134").Value = 2071.484
$ws.Range("I134").Value = 1957.2963
$ws.Range("J134").Value = 2842.25
$ws.Range("K134").Value = 5871.8889
$ws.Range("L134").Value = 8526.75
$ws.Range("M134").Value = -3336.8889
$ws.Range("N134").Value = -13596.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 546.2727
$ws.Range("I105").Value = 439.1875
$ws.Range("J105").Value = 831.8333
$ws.Range("K105").Value = 439.1875
$ws.Range("L105").Value = 831.8333
$ws.Range("M105").Value = 1307.8125
$ws.Range("N105").Value = -4325.8333

$ws.Range("H132").Value = 3751.2163
$ws.Range("I132").Value = 4207.8887
$ws.Range("J132").Value = 2518.2
$ws.Range("K132").Value = 12623.6661
$ws.Range("L132").Value = 7554.599999999999
$ws.Range("M132").Value = -10093.6661
$ws.Range("N132").Value = -12614.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 83333464
$ws.Range("I23").Value = 106.666664
$ws.Range("J23").Value = 166666820
$ws.Range("K23").Value = 319.999992
$ws.Range("L23").Value = 500000460
$ws.Range("M23").Value = -84.99999200000002
$ws.Range("N23").Value = -500000930

$ws.Range("H113").Value = 2463490.2
$ws.Range("I113").Value = 5747460
$ws.Range("J113").Value = 513
$ws.Range("K113").Value = 17242380
$ws.Range("L113").Value = 1539
$ws.Range("M113").Value = -17240210
$ws.Range("N113").Value = -5879

$ws.Range("H122").Value = 972442.7
$ws.Range("I122").Value = 6571.4443
$ws.Range("J122").Value = 4133475.8
$ws.Range("K122").Value = 59142.9987
$ws.Range("L122").Value = 37201282.2
$ws.Range("M122").Value = -56692.9987
$ws.Range("N122").Value = -37206182.2

$ws.Range("H129").Value = 2197.2415
$ws.Range("I129").Value = 1588.3334
$ws.Range("J129").Value = 2356.087
$ws.Range("K129").Value = 4765.0002
$ws.Range("L129").Value = 7068.261
$ws.Range("M129").Value = 234.9997999999996
$ws.Range("N129").Value = -17068.261

$ws.Range("H131").Value = 857.3509
$ws.Range("I131").Value = 403.33334
$ws.Range("J131").Value = 910.7646999999999
$ws.Range("K131").Value = 1210.00002
$ws.Range("L131").Value = 2732.2941
$ws.Range("M131").Value = 3829.99998
$ws.Range("N131").Value = -12812.2941

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 168.75
$ws.Range("I107").Value = 145
$ws.Range("J107").Value = 287.5
$ws.Range("K107").Value = 145
$ws.Range("L107").Value = 287.5
$ws.Range("M107").Value = 1775
$ws.Range("N107").Value = -4127.5

$ws.Range("H122").Value = 4170899.5
$ws.Range("I122").Value = 10001481
$ws.Range("J122").Value = 6198.4287
$ws.Range("K122").Value = 30004443
$ws.Range("L122").Value = 18595.2861
$ws.Range("M122").Value = -30001993
$ws.Range("N122").Value = -23495.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1866.6666
$ws.Range("I61").Value = 1866.6666
$ws.Range("K61").Value = 1866.6666
$ws.Range("M61").Value = -1664.6666

$ws.Range("H113").Value = 1866.6666
$ws.Range("I113").Value = 1866.6666
$ws.Range("K113").Value = 1866.6666
$ws.Range("M113").Value = 303.3334

$ws.Range("H133").Value = 48024.234
$ws.Range("J133").Value = 48024.234
$ws.Range("L133").Value = 48024.234
$ws.Range("N133").Value = -53084.234

$ws.Range("H139").Value = 20218.678
$ws.Range("J139").Value = 20218.678
$ws.Range("L139").Value = 20218.678
$ws.Range("N139").Value = -30498.678

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 281.7857
$ws.Range("I107").Value = 239.4
$ws.Range("J107").Value = 387.75
$ws.Range("K107").Value = 718.2
$ws.Range("L107").Value = 1163.25
$ws.Range("M107").Value = 1201.8
$ws.Range("N107").Value = -5003.25

$ws.Range("H122").Value = 1900.4584
$ws.Range("I122").Value = 1829
$ws.Range("J122").Value = 2000.5
$ws.Range("K122").Value = 5487
$ws.Range("L122").Value = 6001.5
$ws.Range("M122").Value = -3037
$ws.Range("N122").Value = -10901.5

$ws.Range("H123").Value = 23924.092
$ws.Range("I123").Value = 20000
$ws.Range("J123").Value = 24796.111
$ws.Range("K123").Value = 20000
$ws.Range("L123").Value = 24796.111
$ws.Range("M123").Value = -15100
$ws.Range("N123").Value = -34596.111

$ws.Range("H132").Value = 1621.4615
$ws.Range("I132").Value = 945.14703
$ws.Range("J132").Value = 2898.9443
$ws.Range("K132").Value = 2835.44109
$ws.Range("L132").Value = 8696.832900000001
$ws.Range("M132").Value = -305.4410899999998
$ws.Range("N132").Value = -13756.8329

$ws.Range("H136").Value = 4013.4146
$ws.Range("I136").Value = 4501.2856
$ws.Range("J136").Value = 2962.6155
$ws.Range("K136").Value = 13503.8568
$ws.Range("L136").Value = 8887.8465
$ws.Range("M136").Value = -10953.8568
$ws.Range("N136").Value = -13987.8465

$ws.Range("H139").Value = 6866.375
$ws.Range("J139").Value = 6866.375
$ws.Range("L139").Value = 6866.375
